{"js": "const replacements = [\n  [\"2025-12-26 Friday\", \"2025-12-27 Saturday\"],\n  [\"306\u00d79=2754\", \"752\u00d77=5264\"],\n  [\"516\u00d75=2580\", \"988\u00d78=7904\"],\n  [\"142\u00d72=284\", \"444\u00d78=3552\"],\n  [\"601\u00d75=3005\", \"892\u00d75=4460\"],\n  [\"835\u00d72=1670\", \"998\u00d74=3992\"],\n  [\"529\u00d78=4232\", \"389\u00d76=2334\"],\n  [\"305\u00d76=1830\", \"278\u00d77=1946\"],\n  [\"235\u00d74=940\", \"943\u00d78=7544\"],\n  [\"620\u00d77=4340\", \"494\u00d72=988\"],\n  [\"238\u00d79=2142\", \"432\u00d79=3888\"],\n  [\"838\u00d73=2514\", \"260\u00d72=520\"],\n  [\"834\u00d78=6672\", \"788\u00d72=1576\"],\n  [\"723\u00d79=6507\", \"466\u00d76=2796\"],\n  [\"845\u00d73=2535\", \"643\u00d78=5144\"],\n  [\"338\u00d78=2704\", \"369\u00d74=1476\"],\n  [\"231\u00d73=693\", \"866\u00d76=5196\"],\n  [\"279\u00d73=837\", \"895\u00d78=7160\"],\n  [\"906\u00d75=4530\", \"983\u00d75=4915\"],\n  [\"483\u00d76=2898\", \"335\u00d76=2010\"],\n  [\"878\u00d76=5268\", \"786\u00d74=3144\"],\n  [\"849\u00d77=5943\", \"878\u00d79=7902\"],\n  [\"868\u00d77=6076\", \"582\u00d78=4656\"],\n  [\"507\u00d74=2028\", \"480\u00d72=960\"],\n  [\"574\u00d73=1722\", \"246\u00d78=1968\"],\n  [\"955\u00d78=7640\", \"274\u00d79=2466\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-12-26 Friday\", \"2025-12-27 Saturday\"),\n    @(\"306\u00d79=2754\", \"752\u00d77=5264\"),\n    @(\"516\u00d75=2580\", \"988\u00d78=7904\"),\n    @(\"142\u00d72=284\", \"444\u00d78=3552\"),\n    @(\"601\u00d75=3005\", \"892\u00d75=4460\"),\n    @(\"835\u00d72=1670\", \"998\u00d74=3992\"),\n    @(\"529\u00d78=4232\", \"389\u00d76=2334\"),\n    @(\"305\u00d76=1830\", \"278\u00d77=1946\"),\n    @(\"235\u00d74=940\", \"943\u00d78=7544\"),\n    @(\"620\u00d77=4340\", \"494\u00d72=988\"),\n    @(\"238\u00d79=2142\", \"432\u00d79=3888\"),\n    @(\"838\u00d73=2514\", \"260\u00d72=520\"),\n    @(\"834\u00d78=6672\", \"788\u00d72=1576\"),\n    @(\"723\u00d79=6507\", \"466\u00d76=2796\"),\n    @(\"845\u00d73=2535\", \"643\u00d78=5144\"),\n    @(\"338\u00d78=2704\", \"369\u00d74=1476\"),\n    @(\"231\u00d73=693\", \"866\u00d76=5196\"),\n    @(\"279\u00d73=837\", \"895\u00d78=7160\"),\n    @(\"906\u00d75=4530\", \"983\u00d75=4915\"),\n    @(\"483\u00d76=2898\", \"335\u00d76=2010\"),\n    @(\"878\u00d76=5268\", \"786\u00d74=3144\"),\n    @(\"849\u00d77=5943\", \"878\u00d79=7902\"),\n    @(\"868\u00d77=6076\", \"582\u00d78=4656\"),\n    @(\"507\u00d74=2028\", \"480\u00d72=960\"),\n    @(\"574\u00d73=1722\", \"246\u00d78=1968\"),\n    @(\"955\u00d78=7640\", \"274\u00d79=2466\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 0, $false, $find.Replacement.Text, 2)\n}"}
